$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-word the two existing Energizer cell rows to use commas as separators
$ws.Range("A30").Value = "CGAS007 Energizer Single Cell, 3.7 V, Lithium Ion, 1 Ah"
$ws.Range("A31").Value = "CA5L Energizer Single Cell, 3.7 V, Lithium Ion, 980 mAh"

# Prime rows 32:34 with the same formatting (vertical-centre + wrap text)
# as the rest of the table by copying the format from the row above
$ws.Range("A31:E31").Copy() | Out-Null
$ws.Range("A32:E34").PasteSpecial(-4122) | Out-Null

# New row: 2x Clip on Pod Enclosures
$ws.Range("A32").Value = "2x Clip on Pod Enclosures"
$ws.Range("B32").Value = "ENC_001"
$ws.Range("C32").Value = 2.8544
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 1

# New row: 2x Insole with Removable Pod Enclosure
$ws.Range("A33").Value = "2x Insole with Removable Pod Enclosure"
$ws.Range("B33").Value = "ENC_002"
$ws.Range("C33").Value = 8.0434000000000001
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 1

# New row: 2x Insole with Attachable Pod Enclosure
$ws.Range("A34").Value = "2x Insole with Attachable Pod Enclosure"
$ws.Range("B34").Value = "ENC_003"
$ws.Range("C34").Value = 7.5095999999999998
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 1

# Update the sheet view: scroll back to the top-left and select H29
$ws.Activate() | Out-Null
$ws.Range("H29").Select() | Out-Null
